# "update pelajaran web design fundamental semester 1, pertemuan kedua"
# - table
# - pembuatan kolom menggunakan div
# - styling margin
# - github
#
# Adds a new "Sheet2" (group member roster) and extends the "HTML TAG"
# cheatsheet on Sheet1 with four new tag entries: <div>, <span>, <table>,
# <tr>, <td>.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Add Sheet2 (group members) after the existing sheet.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Ricky"
$ws2.Range("A2").Value = "Titan"

# ---------------------------------------------------------------------
# 2) Extend Sheet1 ("HTML TAG") with new rows 17-31, mirroring the
#    existing bordered-table block formatting (rows 6-16).
# ---------------------------------------------------------------------

# Source cells carrying the three alternating border styles already used
# throughout the sheet's cheatsheet block.
$bTop    = "B10"   # blank separator row, top style
$cTop    = "C10"
$bLabel  = "B6"    # tag-name cell style
$cLabel  = "C6"    # description cell style
$bBlank  = "B9"    # blank continuation row, bottom style
$cBlank  = "C9"

function Format-Row($row, $bFmt, $cFmt) {
    $b = $ws1.Range("B" + $row)
    $c = $ws1.Range("C" + $row)
    $ws1.Range($bFmt).Copy()
    $b.PasteSpecial(-4122)   # xlPasteFormats
    $ws1.Range($cFmt).Copy()
    $c.PasteSpecial(-4122)   # xlPasteFormats
}

# -- <div> block (rows 17-19) ------------------------------------------------
Format-Row 17 $bTop   $cTop
Format-Row 18 $bLabel $cLabel
Format-Row 19 $bBlank $cBlank

$ws1.Range("B18").Value = "<div>"
$ws1.Range("C18").Value = "Reserve display website, full, dari kiri hingga kanan website"

# -- <span> block (rows 20-22) ------------------------------------------------
Format-Row 20 $bTop   $cTop
Format-Row 21 $bLabel $cLabel
Format-Row 22 $bBlank $cBlank

$ws1.Range("B21").Value = "<span>"

# ---------------------------------------------------------------------
# 3) Sheet2 members keep being filled in (interleaved in the original
#    editing session, reflected in shared-string ordering).
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = "Dylan"
$ws2.Range("A4").Value = "Andrea"

$ws1.Range("C21").Value = "Reserve display website, hanya sepanjang content saja."
$ws1.Range("C22").Value = "Keterangan: inline"

# -- <table> block (rows 23-25) ------------------------------------------------
Format-Row 23 $bTop   $cTop
Format-Row 24 $bLabel $cLabel
Format-Row 25 $bBlank $cBlank

$ws1.Range("B24").Value = "<table>"
$ws1.Range("C24").Value = "Untuk membuat tabel."

# -- <tr> block (rows 26-28) ------------------------------------------------
Format-Row 26 $bTop   $cTop
Format-Row 27 $bLabel $cLabel
Format-Row 28 $bBlank $cBlank

$ws1.Range("B27").Value = "<tr>"
$ws1.Range("C27").Value = "Untuk membuat row baru dalam 1 table"

# -- <td> block (rows 29-31) ------------------------------------------------
Format-Row 29 $bTop   $cTop
Format-Row 30 $bLabel $cLabel
Format-Row 31 $bBlank $cBlank

$ws1.Range("B30").Value = "<td>"
$ws1.Range("C30").Value = "Untuk membuat cell baru dalam 1 row di dalam 1 table"

# ---------------------------------------------------------------------
# 4) View state: scroll sheet1 to the new rows and zoom in, then restore
#    selection roughly matching the saved workbook state.
# ---------------------------------------------------------------------
$ws1.Range("A21").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150
$ws1.Range("C4:C31").Select() | Out-Null

$ws2.Range("A5").Select() | Out-Null

$ws1.Select() | Out-Null
